# The HR Management application allows HR to log in with a valid username
# and password. This adds a second worksheet ("InvalidUsers") next to the
# existing "Users" sheet, seeded with a set of invalid/unauthorized
# username+password combinations for negative-path testing.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after "Users" so the tab order matches
# (Users, InvalidUsers).
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "InvalidUsers"

# Match the outline defaults used on the "Users" sheet.
$newSheet.Outline.SummaryRow = -4118    # xlSummaryAbove -> summaryBelow="0"
$newSheet.Outline.SummaryColumn = -4131 # xlSummaryOnLeft -> summaryRight="0"

# Header row
$newSheet.Range("A1").Value = "userName"
$newSheet.Range("B1").Value = "Password"

# Data rows
$newSheet.Range("A2").Value = "admin"
$newSheet.Range("B2").Value = "amod12"

$newSheet.Range("A3").Value = "pathi123"
$newSheet.Range("B3").Value = "manager"

$newSheet.Range("A4").Value = "lakshman"
$newSheet.Range("B4").Value = "hulzi14"

# Reuse the same cell formatting/style as the "Users" sheet.
$ws1.Range("A1").Copy() | Out-Null
$newSheet.Range("A1:B4").PasteSpecial(-4122) | Out-Null

# Restore the selections so the saved view matches the authored workbook:
# "Users" still has A3:XFD3 selected, while "InvalidUsers" (now the active
# tab) has B4 selected.
$ws1.Range("A3:XFD3").Select() | Out-Null
$newSheet.Range("B4").Select() | Out-Null
